$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 277
$ws.Range("B277").Value = 61610
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18

# Row 278
$ws.Range("B278").Value = 63565
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6

# Row 292
$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 69
$ws.Range("G292").Value = 9955.32

# Row 293
$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32

# Row 294
$ws.Range("B294").Value = 57802
$ws.Range("E294").Value = 162.71
$ws.Range("F294").Value = -79
$ws.Range("G294").Value = -11334.92

# Row 295
$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 0
$ws.Range("G295").Value = 0

# Row 296
$ws.Range("B296").Value = 63531
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 79
$ws.Range("G296").Value = 11334.92

# Row 299
$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 132
$ws.Range("G299").Value = 6288.48

# Row 300
$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12

# Row 311
$ws.Range("B311").Value = 63563
$ws.Range("E311").Value = 119.04
$ws.Range("F311").Value = 0
$ws.Range("G311").Value = 0

# Row 312
$ws.Range("B312").Value = 61605
$ws.Range("E312").Value = 133.78
$ws.Range("F312").Value = -13
$ws.Range("G312").Value = -1455.48

# Row 420
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 41
$ws.Range("G420").Value = 4327.14

# Row 421
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2

# Row 465
$ws.Range("B465").Value = 65069
$ws.Range("E465").Value = 14.3
$ws.Range("F465").Value = 2
$ws.Range("G465").Value = 26.9

# Row 466
$ws.Range("B466").Value = 53757
$ws.Range("E466").Value = 16.08
$ws.Range("F466").Value = -159
$ws.Range("G466").Value = -2138.55

# Row 467
$ws.Range("B467").Value = 65068
$ws.Range("E467").Value = 13.97
$ws.Range("F467").Value = 63
$ws.Range("G467").Value = 828.45

# Row 468
$ws.Range("B468").Value = 53602
$ws.Range("E468").Value = 15.69
$ws.Range("F468").Value = -231
$ws.Range("G468").Value = -3037.65

# Row 472
$ws.Range("B472").Value = 64915
$ws.Range("E472").Value = 20.98
$ws.Range("F472").Value = 0
$ws.Range("G472").Value = 0

# Row 473
$ws.Range("B473").Value = 45695
$ws.Range("E473").Value = 23.58
$ws.Range("F473").Value = -36
$ws.Range("G473").Value = -710.28

# Row 479
$ws.Range("B479").Value = 45718
$ws.Range("E479").Value = 19.38
$ws.Range("F479").Value = -294
$ws.Range("G479").Value = -4768.68

# Row 480
$ws.Range("B480").Value = 64927
$ws.Range("E480").Value = 17.26
$ws.Range("F480").Value = 119
$ws.Range("G480").Value = 1930.18

# Row 490
$ws.Range("B490").Value = 65067
$ws.Range("E490").Value = 15.65
$ws.Range("F490").Value = 172
$ws.Range("G490").Value = 2533.56

# Row 491
$ws.Range("B491").Value = 53595
$ws.Range("E491").Value = 17.61
$ws.Range("F491").Value = -335
$ws.Range("G491").Value = -4934.55

# Row 595
$ws.Range("B595").Value = 64836
$ws.Range("E595").Value = 104.71
$ws.Range("F595").Value = 0
$ws.Range("G595").Value = 0

# Row 596
$ws.Range("B596").Value = 60031
$ws.Range("E596").Value = 111.69
$ws.Range("F596").Value = -5
$ws.Range("G596").Value = -492.5

# Row 600
$ws.Range("B600").Value = 60022
$ws.Range("E600").Value = 37.22
$ws.Range("F600").Value = -113
$ws.Range("G600").Value = -3709.79

# Row 601
$ws.Range("B601").Value = 64830
$ws.Range("E601").Value = 34.9
$ws.Range("F601").Value = 111
$ws.Range("G601").Value = 3644.13

# Row 736
$ws.Range("B736").Value = 65362
$ws.Range("F736").Value = 37
$ws.Range("G736").Value = 1512.19

# Row 737
$ws.Range("B737").Value = 65079
$ws.Range("F737").Value = 21
$ws.Range("G737").Value = 858.27

